# Action List workbook update:
#  - Add a new "Logon" sheet (after "Home Page") documenting logon testing
#  - Add a "Basic functionality" note on the Home Page sheet
#  - Add a "Menu not working" note on the meeting space sheet

$wb = $excel.ActiveWorkbook

$home    = $wb.Worksheets.Item("Home Page")
$meeting = $wb.Worksheets.Item("meeting space")

# ---------------------------------------------------------------------------
# Insert new "Logon" worksheet right after "Home Page"
# ---------------------------------------------------------------------------
$logon = $wb.Worksheets.Add($null, $home)
$logon.Name = "Logon"

$logon.Columns("C").ColumnWidth = 49.736979166666664
$logon.Columns("D").ColumnWidth = 40.166666666666664
$logon.Columns("E").ColumnWidth = 39.166666666666664
$logon.Columns("E").WrapText = $true

# ---------------------------------------------------------------------------
# Fill in the cell content in the same order it was originally authored, so
# that new shared-string entries come out in the expected order.
# ---------------------------------------------------------------------------
$meeting.Range("B10").Value = "Menu not working "

$logon.Range("D3").Value = "Forgot Password Option "
$logon.Range("E3").Value = "Allowed me to enter email but did not receive email"

$home.Range("B7").Value = "Basic functionality of this Page only"
$home.Range("B7").Interior.Color = 65535

$logon.Range("C3").Value = " Logon  Type email Deirdre.Shanahan@gmail.com"
$logon.Range("B3").Value = "Test1"
$logon.Range("B6").Value = "Test2"

$logon.Range("D1").Value = "Option"
$logon.Range("E1").Value = "Comments"
$logon.Range("F1").Value = "Action"

$logon.Range("C6").Value = "Logon as Lisa Sharp Password welcomE01#"
$logon.Range("D6").Value = "allowed me to logon and took me to Meeting Space "
$logon.Range("D7").Value = "Welcome Lisa at the top of meeting space"
$logon.Range("E6").Value = "looks good so far"

$logon.Range("B9").Value = "Test 3"
$logon.Range("C9").Value = "Logon as  Carol.James@gmail.com password welcomE01#"
$logon.Range("D9").Value = "Logon opttion"
$logon.Range("E9").Value = "Brought me to meeting space but error because of Profile Screen 2"

# Cells re-using already existing shared strings ("Open" / " ")
$logon.Range("F3").Value = "Open"
$logon.Range("C4").Value = " "
$logon.Range("D4").Value = " "
$logon.Range("E4").Value = " "
$logon.Range("F4").Value = " "
$logon.Range("F9").Value = "Open"
$meeting.Range("C10").Value = "Open"

$logon.Rows(9).RowHeight = 30
$logon.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selections / active sheet to match the saved view state
# ---------------------------------------------------------------------------
$home.Range("B7").Select() | Out-Null
$meeting.Range("C11").Select() | Out-Null

$logon.Activate() | Out-Null
$logon.Range("E9").Select() | Out-Null
